$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'sliding pads for exercise'
    2 = 'compression capri men'
    3 = 'basketball pants for boys'
    4 = 'padded baseball pants'
    5 = 'cycling capri pants'
    6 = 'sliding pants baseball mens'
    7 = 'below the knee shorts for men'
    8 = 'mens long basketball shorts below knee'
    9 = 'youth small knee pads basketball'
    10 = 'basketball leg pads'
    11 = 'short baseball pants'
    12 = 'basketball bump'
    13 = 'knee pads running'
    14 = 'football compression shorts with pads'
    15 = 'knee high baseball pants mens'
    16 = 'knee padded compression'
    17 = 'athletic pads'
    18 = 'mens protection pads'
    19 = 'youth knee pads wrestling'
    20 = 'men softball pants'
    21 = 'baseball sliding pants'
    22 = 'boys sliding pants'
    23 = 'youth boys basketball pants'
    24 = 'lacrosse sweat pants'
    25 = 'raceface knee pads'
    26 = 'capri compression pants men'
    27 = 'snowboarding compression pants'
    28 = 'men knee pad pants'
    29 = 'youth basketball tights for boys'
    30 = 'knee pads for basketball youth'
    31 = 'mtn bike knee pads'
    32 = 'baseball padded sliding shorts'
    33 = 'leggings with baseballs'
    34 = 'youth knee pad wrestling'
    35 = '28 basketball'
    36 = 'tights with knee'
    37 = 'compression capris'
    38 = 'padded tights for football'
    39 = 'baseball tights for boys'
    40 = 'baseball sliding shorts men'
    41 = 'knee pads for basketball youth boys'
    42 = 'youth leggings boys basketball'
    43 = 'basketballs leggings'
    44 = 'softball items'
    45 = 'football leggings for men'
    46 = 'basketball compression pants youth'
    47 = 'compression pants men black'
    48 = 'mens work pants with knee pads'
    49 = 'youth tights'
    50 = 'mens down pants'
    51 = 'compression basketball pants youth'
    52 = 'basketball compression pants women'
    53 = 'softball slider'
    54 = 'work knee pads under pants'
    55 = 'calf compression pants'
    56 = 'leg pads basketball'
    57 = 'compression pants men football'
    58 = 'softball shorts men'
    59 = 'mens softball shorts'
    60 = 'mens capri shorts below knee'
    61 = 'compression tights youth'
    62 = 'adidas knee pads'
    63 = 'athletic leggings mens'
    64 = 'knee compression shorts'
    65 = 'compression knee pads pair'
    66 = 'indoor volleyball knee pads'
    67 = 'youth basketball pants boys'
    68 = 'softball sliding shorts girls padded'
    69 = 'lacrosse compression shorts padded'
    70 = 'taken leggings'
    71 = 'mens basketball tights'
    72 = 'black compression pants men'
    73 = 'softball sliding shorts'
    74 = 'girls basketball knee pads youth'
    75 = 'kneepads basketball'
    76 = 'knee protector for construction'
    77 = 'knee pad sleeve basketball'
    78 = 'youth basketball leggings'
    79 = 'medium compression pants'
    80 = 'football knee pads for men'
    81 = 'mens compression leggings'
    82 = 'knee pad for yoga'
    83 = 'boys knee pads basketball'
    84 = 'baseball sliding'
    85 = 'knee pads for biking men'
    86 = 'knee pads girls basketball'
    87 = 'black football pants'
    88 = 'lacrosse pants'
    89 = 'team work softball pants'
    90 = 'long basketball shorts for men below knee'
    91 = 'basketball hex pads'
    92 = 'compression shorts padded basketball'
    93 = 'big boys tights'
    94 = 'basketball tights'
    95 = 'hockey knee pads adult'
    96 = 'padded compression pants football'
    97 = 'hockey hip pads adult'
    98 = 'soccer goalkeeper pads'
    99 = 'athletic mens leggings'
    100 = 'knee pads flexible'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
